# Generate Report for Handoff
# Update the "Latest Handoff Datetime" (column D) timestamps for the files that
# were just (re-)handed off, on both the zh-cn and de-de status sheets.
# Rows 4,6,7,8,9,10 correspond to the files being handed off in this report run
# (row 5 - "In Translation" - is intentionally left untouched).

$wb = $excel.ActiveWorkbook

$rowsToUpdate = @(4, 6, 7, 8, 9, 10)

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rowsToUpdate) {
    $ws_zhcn.Cells.Item($r, 4).Value = "2016-03-04 04:22:08"
}

$ws_dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rowsToUpdate) {
    $ws_dede.Cells.Item($r, 4).Value = "2016-03-04 04:22:23"
}
